# Weekly update: add a new week's price record (row 19) for
# Hortaliza, Vega Monumental Concepción - Poroto verde.
# All existing data rows from 19 downward shift down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19, pushing existing rows 19-27 to 20-28.
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Range("A19").Value = 11
$ws.Range("B19").Value = 'Vega Monumental Concepción'
$ws.Range("C19").Value = 'Bíobío'
$ws.Range("D19").Value = 44489
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = 'Poroto verde'
$ws.Range("H19").Value = 'Magnum'
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 40000
$ws.Range("L19").Value = 42000
$ws.Range("M19").Value = 41000
$ws.Range("N19").Value = '$/malla 25 kilos'
$ws.Range("O19").Value = 'Perú'
$ws.Range("P19").Value = 1640
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = 'Hortaliza'
